# Auto-generated Excel COM-interop edit script
# Applies numeric cell updates (and a few cell adds/clears) across all 8 job sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 598671.9399999999
$ws.Range("J17").Value = 662601.2
$ws.Range("L17").Value = 1987803.6
$ws.Range("N17").Value = -1988139.6

$ws.Range("H62").Value = 9926.933999999999
$ws.Range("I62").Value = 9863
$ws.Range("K62").Value = 9863
$ws.Range("M62").Value = -9239

$ws.Range("H65").Value = 9926.933999999999
$ws.Range("I65").Value = 9863
$ws.Range("K65").Value = 49315
$ws.Range("M65").Value = -46195

$ws.Range("H129").Value = 3782.611
$ws.Range("I129").Value = 2258
$ws.Range("K129").Value = 6774
$ws.Range("M129").Value = -1774

$ws.Range("H137").Value = 9095.286
$ws.Range("I137").Value = 2892.3333
$ws.Range("J137").Value = 10787
$ws.Range("K137").Value = 8676.999899999999
$ws.Range("L137").Value = 32361
$ws.Range("M137").Value = -6126.999899999999
$ws.Range("N137").Value = -37461

$ws.Range("H138").Value = 3122.617
$ws.Range("I138").Value = 1397.375
$ws.Range("J138").Value = 3714.1287
$ws.Range("K138").Value = 4192.125
$ws.Range("L138").Value = 11142.3861
$ws.Range("M138").Value = 947.875
$ws.Range("N138").Value = -21422.3861


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1657.9259
$ws.Range("I2").Value = 1782.2632
$ws.Range("K2").Value = 1782.2632
$ws.Range("M2").Value = -1669.2632

$ws.Range("H32").Value = 3381.756
$ws.Range("I32").Value = 2947.4443
$ws.Range("K32").Value = 2947.4443
$ws.Range("M32").Value = -2660.4443

$ws.Range("H45").Value = 15444.389
$ws.Range("I45").Value = 11499
$ws.Range("K45").Value = 11499
$ws.Range("M45").Value = -11122

$ws.Range("H97").Value = 1139.5358
$ws.Range("I97").Value = 1156.5416
$ws.Range("K97").Value = 1156.5416
$ws.Range("M97").Value = -660.5416

$ws.Range("H116").Value = 1657.9259
$ws.Range("I116").Value = 1782.2632
$ws.Range("K116").Value = 1782.2632
$ws.Range("M116").Value = 511.7367999999999

$ws.Range("H122").Value = 3031.3
$ws.Range("I122").Value = 3112.375
$ws.Range("K122").Value = 9337.125
$ws.Range("M122").Value = -6887.125


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1657.9259
$ws.Range("I3").Value = 1782.2632
$ws.Range("K3").Value = 1782.2632
$ws.Range("M3").Value = -1668.2632

$ws.Range("H20").Value = 2211.9795
$ws.Range("I20").Value = 1889.8572
$ws.Range("K20").Value = 1889.8572
$ws.Range("M20").Value = -1642.8572

$ws.Range("H86").Value = 7828.4634
$ws.Range("I86").Value = 5824.593
$ws.Range("K86").Value = 5824.593
$ws.Range("M86").Value = -4701.593

$ws.Range("H89").Value = 7828.4634
$ws.Range("I89").Value = 5824.593
$ws.Range("K89").Value = 29122.965
$ws.Range("M89").Value = -23506.965


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 5000
$ws.Range("I58").Value = 5000
$ws.Range("K58").Value = 5000
$ws.Range("M58").Value = -4797

$ws.Range("H99").Value = 3799.125
$ws.Range("I99").Value = 3799.125
$ws.Range("K99").Value = 3799.125
$ws.Range("M99").Value = -2301.125

$ws.Range("H126").Value = 3799.125
$ws.Range("I126").Value = 3799.125
$ws.Range("K126").Value = 11397.375
$ws.Range("M126").Value = -8927.375

$ws.Range("H132").Value = 3676.8572
$ws.Range("I132").Value = 3373.0833
$ws.Range("K132").Value = 10119.2499
$ws.Range("M132").Value = -7589.249899999999

$ws.Range("H134").Value = 10356.742
$ws.Range("I134").Value = 2335.7917
$ws.Range("K134").Value = 7007.375100000001
$ws.Range("M134").Value = -4472.375100000001

$ws.Range("H136").Value = 5000
$ws.Range("I136").Value = 5000
$ws.Range("K136").Value = 15000
$ws.Range("M136").Value = -12450


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H42").Value = 8951.5
$ws.Range("J42").Value = 15000
$ws.Range("L42").Value = 45000
$ws.Range("N42").Value = -46068

$ws.Range("H82").Value = 17240.4
$ws.Range("I82").Value = 17240.4
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 51721.2
$ws.Range("L82").Value = 0
$ws.Range("M82").ClearContents() | Out-Null
$ws.Range("N82").Value = -51315.2

$ws.Range("H85").Value = 17240.4
$ws.Range("I85").Value = 17240.4
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 51721.2
$ws.Range("L85").Value = 0
$ws.Range("M85").ClearContents() | Out-Null
$ws.Range("N85").Value = -50317.2

$ws.Range("H112").Value = 9003.25
$ws.Range("I112").Value = 8999.5
$ws.Range("J112").Value = 9004.5
$ws.Range("K112").Value = 26998.5
$ws.Range("L112").Value = 27013.5
$ws.Range("M112").Value = -25890.5
$ws.Range("N112").Value = -29229.5

$ws.Range("H121").Value = 2935.158
$ws.Range("I121").Value = 378.5
$ws.Range("J121").Value = 3616.9333
$ws.Range("K121").Value = 1135.5
$ws.Range("L121").Value = 10850.7999
$ws.Range("M121").Value = 174.5
$ws.Range("N121").Value = -13470.7999

$ws.Range("H141").Value = 146288.72
$ws.Range("I141").Value = 2996
$ws.Range("J141").Value = 253758.25
$ws.Range("K141").Value = 8988
$ws.Range("L141").Value = 761274.75
$ws.Range("M141").Value = -3808
$ws.Range("N141").Value = -771634.75


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7063.316
$ws.Range("I70").Value = 7434.3335
$ws.Range("J70").Value = 6729.4
$ws.Range("K70").Value = 7434.3335
$ws.Range("L70").Value = 6729.4
$ws.Range("M70").Value = -7164.3335
$ws.Range("N70").Value = -7269.4

$ws.Range("H73").Value = 7063.316
$ws.Range("I73").Value = 7434.3335
$ws.Range("J73").Value = 6729.4
$ws.Range("K73").Value = 7434.3335
$ws.Range("L73").Value = 6729.4
$ws.Range("M73").Value = -6498.3335
$ws.Range("N73").Value = -8601.4

$ws.Range("H132").Value = 180899.02
$ws.Range("I132").Value = 190827.28
$ws.Range("K132").Value = 572481.84
$ws.Range("M132").Value = -569951.84


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 17753.223
$ws.Range("I7").Value = 17753.223
$ws.Range("K7").Value = 17753.223
$ws.Range("M7").Value = -17641.223

$ws.Range("H40").Value = 5996.5713
$ws.Range("I40").Value = 5996.5713
$ws.Range("K40").Value = 5996.5713
$ws.Range("M40").Value = -5860.5713

$ws.Range("H46").Value = 3979.9
$ws.Range("I46").Value = 1828.5714
$ws.Range("J46").Value = 8999.666999999999
$ws.Range("K46").Value = 1828.5714
$ws.Range("L46").Value = 8999.666999999999
$ws.Range("M46").Value = -1640.5714
$ws.Range("N46").Value = -9375.666999999999

$ws.Range("H55").Value = 182.2
$ws.Range("J55").Value = 400
$ws.Range("L55").Value = 400
$ws.Range("N55").Value = -746

$ws.Range("H61").Value = 3700.875
$ws.Range("I61").Value = 3518.111
$ws.Range("K61").Value = 3518.111
$ws.Range("M61").Value = -3316.111

$ws.Range("H74").Value = 33711.8
$ws.Range("I74").Value = 33523
$ws.Range("J74").Value = 33995
$ws.Range("K74").Value = 33523
$ws.Range("L74").Value = 33995
$ws.Range("M74").Value = -32525
$ws.Range("N74").Value = -35991

$ws.Range("H77").Value = 33711.8
$ws.Range("I77").Value = 33523
$ws.Range("J77").Value = 33995
$ws.Range("K77").Value = 100569
$ws.Range("L77").Value = 101985
$ws.Range("M77").Value = -95577
$ws.Range("N77").Value = -111969

$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").ClearContents() | Out-Null
$ws.Range("N80").Value = 0

$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").ClearContents() | Out-Null
$ws.Range("N83").Value = 0

$ws.Range("H93").Value = 797059
$ws.Range("I93").Value = 1113297.6
$ws.Range("J93").Value = 6462.5
$ws.Range("K93").Value = 1113297.6
$ws.Range("L93").Value = 6462.5
$ws.Range("M93").Value = -1112049.6
$ws.Range("N93").Value = -8958.5

$ws.Range("H113").Value = 3700.875
$ws.Range("I113").Value = 3518.111
$ws.Range("K113").Value = 3518.111
$ws.Range("M113").Value = -1348.111

$ws.Range("H126").Value = 17753.223
$ws.Range("I126").Value = 17753.223
$ws.Range("K126").Value = 53259.66900000001
$ws.Range("M126").Value = -50789.66900000001

$ws.Range("H136").Value = 7040.7856
$ws.Range("I136").Value = 6797.5835
$ws.Range("K136").Value = 20392.7505
$ws.Range("M136").Value = -17842.7505


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3744.65
$ws.Range("I126").Value = 3478.9333
$ws.Range("K126").Value = 10436.7999
$ws.Range("M126").Value = -7966.7999

